$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns I and J to match the width already used by columns K:N.
$ws.Columns.Item(9).ColumnWidth = 5.7109375
$ws.Columns.Item(10).ColumnWidth = 5.7109375

# Update the row-1 data values.
$ws.Range("B1").Value = 3
$ws.Range("C1").Value = 19
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 31
$ws.Range("F1").Value = 19
$ws.Range("G1").Value = 33
$ws.Range("H1").Value = 31
$ws.Range("I1").Value = 0.08299999999999999
$ws.Range("J1").Value = 0.074999999999999997
$ws.Range("K1").Value = 0.076999999999999985
